$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "61.016.82"
Set-TextValue "E2" "  -0.58%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.390.90"
Set-TextValue "E3" "  -1.34%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "572.01"
Set-TextValue "E5" "  -0.36%  "

# Row 6 - Solana
Set-TextValue "D6" "142.02"
Set-TextValue "E6" "  -1.47%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.391.46"
Set-TextValue "E7" "  -1.38%  "

# Row 9 - XRP
Set-TextValue "E9" "  -0.82%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.62"
Set-TextValue "E10" "  +0.32%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.394"
Set-TextValue "E12" "  +1.33%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.969.77"
Set-TextValue "E13" "  -1.33%  "

# Row 14 - TRON
Set-TextValue "E14" "  +2.31%  "

# Row 15 - Avalanche
Set-TextValue "D15" "27.86"
Set-TextValue "E15" "  -2.29%  "

# Row 16 - was WrappedEther, now ShibaInu
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000170"
Set-TextValue "E16" "  -1.24%  "

# Row 17 - was ShibaInu, now WrappedEther
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.386.94"

# Row 18 - WrappedBTC
Set-TextValue "D18" "61.029.47"
Set-TextValue "E18" "  -0.66%  "

# Row 19 - Polkadot
Set-TextValue "E19" "  -3.52%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.78"
Set-TextValue "E20" "  -4.60%  "

# Row 21 - Uniswap
Set-TextValue "D21" "8.92"
Set-TextValue "E21" "  -4.41%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "382.08"
Set-TextValue "E22" "  -4.74%  "

# Row 23 - Polygon
Set-TextValue "E23" "  -1.28%  "

# Row 24 - Litecoin
Set-TextValue "D24" "74.49"
Set-TextValue "E24" "  +0.60%  "

# Row 25 - Dai
Set-TextValue "E25" "  -0.16%  "

# Row 26 - PEPE
Set-TextValue "E26" "  -4.33%  "

# Row 27 - WrappedeETH
Set-TextValue "D27" "3.531.26"
Set-TextValue "E27" "  -1.50%  "

# Row 28 - Kaspa
Set-TextValue "D28" "0.180"
Set-TextValue "E28" "  +0.70%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "E29" "  -0.04%  "

# Row 30 - RenderToken
Set-TextValue "E30" "  -2.55%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "7.98"
Set-TextValue "E31" "  -2.58%  "

# Row 32 - PancakeSwap
Set-TextValue "E32" "  -0.96%  "

# Row 33 - Fetch.AI
Set-TextValue "E33" "  -7.25%  "

# Row 34 - USDe
Set-TextValue "E34" "  -0.04%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "23.44"
Set-TextValue "E35" "  -2.05%  "

# Row 36 - Aptos
Set-TextValue "E36" "  -0.60%  "

# Row 37 - Monero
Set-TextValue "D37" "167.06"
Set-TextValue "E37" "  +0.37%  "

# Row 38 - RenzoRestakedETH
Set-TextValue "D38" "3.423.25"
Set-TextValue "E38" "  -1.21%  "

# Row 39 - NEARProtocol
Set-TextValue "D39" "5.00"
Set-TextValue "E39" "  -2.49%  "

# Row 40 - ImmutableX
Set-TextValue "E40" "  -4.32%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0769"
Set-TextValue "E41" "  -2.57%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "27.01"
Set-TextValue "E42" "  -0.31%  "

# Row 43 - Mantle
Set-TextValue "E43" "  -2.66%  "

# Row 44 - FirstDigitalUSD
Set-TextValue "E44" "  +0.07%  "

# Row 45 - Filecoin
Set-TextValue "E45" "  -2.09%  "

# Row 46 - Stacks
Set-TextValue "E46" "  -2.93%  "

# Row 47 - ONDO
Set-TextValue "D47" "1.14"
Set-TextValue "E47" "  -1.25%  "

# Row 48 - Maker
Set-TextValue "D48" "2.484.77"
Set-TextValue "E48" "  -4.87%  "

# Row 49 - Cosmos
Set-TextValue "E49" "  -2.16%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "22.86"
Set-TextValue "E50" "  -1.91%  "

# Row 51 - VeChain
Set-TextValue "D51" "0.0265"
Set-TextValue "E51" "  +1.37%  "
